$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("resampling")

# --- Add the new "roberta + lstm" baseline result row ---
$ws.Range("A6").Value = "roberta + lstm"
$ws.Range("B6").Value = 1.4501999999999999
$ws.Range("C6").Value = 0.65380000000000005
$ws.Range("D6").Value = 0.64070000000000005
$ws.Range("E6").Value = 0.65380000000000005
$ws.Range("F6").Value = 0.63419999999999999
$ws.Range("G6").Value = 0.30430000000000001
$ws.Range("H6").Value = 0.31640000000000001
$ws.Range("I6").Value = 0.54590000000000005

# Match the numeric formatting used by the rest of the table (style index 2 -> "0.0000")
$ws.Range("B6:I6").NumberFormat = "0.0000"

# Column A widens to fit the longer label "roberta + lstm" (was sized for "bert + lstm").
# The COM bridge quantizes ColumnWidth onto a 1/7-character pixel grid, so the nearest
# representable value to the authored 14.25 stored width is used here.
$ws.Columns.Item(1).ColumnWidth = 13.5357142857142857

# Selection moves to G11
$ws.Range("G11").Select() | Out-Null
